$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update timestamp footer (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 6 de Mayo de 2020 a las 10:03"

# --- Update country data rows ---
# Row 10: Rusia
$ws.Cells.Item(10, 2).Value = 165929
$ws.Cells.Item(10, 3).Value = 10559
$ws.Cells.Item(10, 4).Value = 21327
$ws.Cells.Item(10, 5).Value = 143065
$ws.Cells.Item(10, 6).Value = 2300
$ws.Cells.Item(10, 7).Value = 86
$ws.Cells.Item(10, 8).Value = 1537

# Row 29: Singapur
$ws.Cells.Item(29, 2).Value = 20198
$ws.Cells.Item(29, 3).Value = 788
$ws.Cells.Item(29, 4).Value = 1519
$ws.Cells.Item(29, 5).Value = 18661
$ws.Cells.Item(29, 6).Value = 24
$ws.Cells.Item(29, 7).Value = 0
$ws.Cells.Item(29, 8).Value = 18

# Row 36: Polonia
$ws.Cells.Item(36, 2).Value = 14431
$ws.Cells.Item(36, 3).Value = 0
$ws.Cells.Item(36, 4).Value = 4655
$ws.Cells.Item(36, 5).Value = 9060
$ws.Cells.Item(36, 6).Value = 160
$ws.Cells.Item(36, 7).Value = 0
$ws.Cells.Item(36, 8).Value = 716

# Row 37: Rumania
$ws.Cells.Item(37, 2).Value = 13837
$ws.Cells.Item(37, 3).Value = 0
$ws.Cells.Item(37, 4).Value = 5454
$ws.Cells.Item(37, 5).Value = 7529
$ws.Cells.Item(37, 6).Value = 244
$ws.Cells.Item(37, 7).Value = 13
$ws.Cells.Item(37, 8).Value = 854

# Row 60: Kazajistan
$ws.Cells.Item(60, 2).Value = 4298
$ws.Cells.Item(60, 3).Value = 93
$ws.Cells.Item(60, 4).Value = 1299
$ws.Cells.Item(60, 5).Value = 2970
$ws.Cells.Item(60, 6).Value = 40
$ws.Cells.Item(60, 7).Value = 0
$ws.Cells.Item(60, 8).Value = 29

# Row 63: Afganistan
$ws.Cells.Item(63, 2).Value = 3392
$ws.Cells.Item(63, 3).Value = 168
$ws.Cells.Item(63, 4).Value = 458
$ws.Cells.Item(63, 5).Value = 2830
$ws.Cells.Item(63, 6).Value = 7
$ws.Cells.Item(63, 7).Value = 9
$ws.Cells.Item(63, 8).Value = 104

# Row 67: Armenia
$ws.Cells.Item(67, 1).Value = "Armenia"
$ws.Cells.Item(67, 2).Value = 2782
$ws.Cells.Item(67, 3).Value = 163
$ws.Cells.Item(67, 4).Value = 1135
$ws.Cells.Item(67, 5).Value = 1607
$ws.Cells.Item(67, 6).Value = 10
$ws.Cells.Item(67, 7).Value = 0
$ws.Cells.Item(67, 8).Value = 40

# Row 68: Oman
$ws.Cells.Item(68, 1).Value = "Oman"
$ws.Cells.Item(68, 2).Value = 2735
$ws.Cells.Item(68, 3).Value = 0
$ws.Cells.Item(68, 4).Value = 858
$ws.Cells.Item(68, 5).Value = 1864
$ws.Cells.Item(68, 6).Value = 17
$ws.Cells.Item(68, 7).Value = 0
$ws.Cells.Item(68, 8).Value = 13

# Row 69: Ghana
$ws.Cells.Item(69, 1).Value = "Ghana"
$ws.Cells.Item(69, 2).Value = 2719
$ws.Cells.Item(69, 3).Value = 0
$ws.Cells.Item(69, 4).Value = 294
$ws.Cells.Item(69, 5).Value = 2407
$ws.Cells.Item(69, 6).Value = 4
$ws.Cells.Item(69, 7).Value = 0
$ws.Cells.Item(69, 8).Value = 18

# Row 70: Grecia
$ws.Cells.Item(70, 1).Value = "Grecia"
$ws.Cells.Item(70, 2).Value = 2642
$ws.Cells.Item(70, 3).Value = 0
$ws.Cells.Item(70, 4).Value = 1374
$ws.Cells.Item(70, 5).Value = 1122
$ws.Cells.Item(70, 6).Value = 35
$ws.Cells.Item(70, 7).Value = 0
$ws.Cells.Item(70, 8).Value = 146

# Row 87: Eslovaquia
$ws.Cells.Item(87, 1).Value = "Eslovaquia"
$ws.Cells.Item(87, 2).Value = 1429
$ws.Cells.Item(87, 3).Value = 8
$ws.Cells.Item(87, 4).Value = 762
$ws.Cells.Item(87, 5).Value = 642
$ws.Cells.Item(87, 6).Value = 4
$ws.Cells.Item(87, 7).Value = 0
$ws.Cells.Item(87, 8).Value = 25

# Row 88: Lituania
$ws.Cells.Item(88, 1).Value = "Lituania"
$ws.Cells.Item(88, 2).Value = 1423
$ws.Cells.Item(88, 3).Value = 0
$ws.Cells.Item(88, 4).Value = 678
$ws.Cells.Item(88, 5).Value = 699
$ws.Cells.Item(88, 6).Value = 17
$ws.Cells.Item(88, 7).Value = 0
$ws.Cells.Item(88, 8).Value = 46

# Row 94: Letonia
$ws.Cells.Item(94, 2).Value = 900
$ws.Cells.Item(94, 3).Value = 4
$ws.Cells.Item(94, 4).Value = 464
$ws.Cells.Item(94, 5).Value = 419
$ws.Cells.Item(94, 6).Value = 3
$ws.Cells.Item(94, 7).Value = 0
$ws.Cells.Item(94, 8).Value = 17

# Row 144: Guadalupe
$ws.Cells.Item(144, 2).Value = 152
$ws.Cells.Item(144, 3).Value = 0
$ws.Cells.Item(144, 4).Value = 104
$ws.Cells.Item(144, 5).Value = 35
$ws.Cells.Item(144, 6).Value = 4
$ws.Cells.Item(144, 7).Value = 1
$ws.Cells.Item(144, 8).Value = 13

# Row 149: Brunei
$ws.Cells.Item(149, 1).Value = "Brunei"
$ws.Cells.Item(149, 2).Value = 139
$ws.Cells.Item(149, 3).Value = 1
$ws.Cells.Item(149, 4).Value = 131
$ws.Cells.Item(149, 5).Value = 7
$ws.Cells.Item(149, 6).Value = 2
$ws.Cells.Item(149, 7).Value = 0
$ws.Cells.Item(149, 8).Value = 1

# Row 150: Guayana Francesa
$ws.Cells.Item(150, 1).Value = "Guayana Francesa"
$ws.Cells.Item(150, 2).Value = 138
$ws.Cells.Item(150, 3).Value = 5
$ws.Cells.Item(150, 4).Value = 111
$ws.Cells.Item(150, 5).Value = 26
$ws.Cells.Item(150, 6).Value = 2
$ws.Cells.Item(150, 7).Value = 0
$ws.Cells.Item(150, 8).Value = 1

# Row 205: Seychelles
$ws.Cells.Item(205, 1).Value = "Seychelles"
$ws.Cells.Item(205, 2).Value = 11
$ws.Cells.Item(205, 3).Value = 0
$ws.Cells.Item(205, 4).Value = 8
$ws.Cells.Item(205, 5).Value = 3
$ws.Cells.Item(205, 6).Value = 0
$ws.Cells.Item(205, 7).Value = 0
$ws.Cells.Item(205, 8).Value = 0

# Row 206: Montserrat
$ws.Cells.Item(206, 1).Value = "Montserrat"
$ws.Cells.Item(206, 2).Value = 11
$ws.Cells.Item(206, 3).Value = 0
$ws.Cells.Item(206, 4).Value = 7
$ws.Cells.Item(206, 5).Value = 3
$ws.Cells.Item(206, 6).Value = 1
$ws.Cells.Item(206, 7).Value = 0
$ws.Cells.Item(206, 8).Value = 1
